$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.171588182449341
$ws.Range("B1").Value = 2.383412837982178
$ws.Range("D1").Value = 2.374018430709839
$ws.Range("E1").Value = 1.210644483566284
